$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For every touched cell, force a text number-format before writing so Excel
# does not auto-coerce numeric-looking strings (e.g. "522.77") into real
# numbers, then restore the default "Normal" style once the value is set.
function Set-TextValue($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws "D2" '58.208.71'
Set-TextValue $ws "E2" '  -0.05%  '
Set-TextValue $ws "D3" '2.595.24'
Set-TextValue $ws "E3" '  -0.52%  '
Set-TextValue $ws "E4" '  +0.09%  '
Set-TextValue $ws "D5" '522.77'
Set-TextValue $ws "E5" '  +0.74%  '
Set-TextValue $ws "E6" '  +0.83%  '
Set-TextValue $ws "E7" '  -0.20%  '
Set-TextValue $ws "D8" '0.569'
Set-TextValue $ws "E8" '  +0.40%  '
Set-TextValue $ws "D9" '2.614.68'
Set-TextValue $ws "E9" '  -0.01%  '
Set-TextValue $ws "D10" '6.63'
Set-TextValue $ws "E10" '  -1.27%  '
Set-TextValue $ws "E11" '  -1.35%  '
Set-TextValue $ws "D12" '0.337'
Set-TextValue $ws "E12" '  -0.23%  '
Set-TextValue $ws "E13" '  -0.29%  '
Set-TextValue $ws "D14" '3.049.22'
Set-TextValue $ws "E14" '  -0.56%  '
Set-TextValue $ws "D15" '58.214.73'
Set-TextValue $ws "E15" '  +0.01%  '
Set-TextValue $ws "D16" '20.53'
Set-TextValue $ws "E16" '  -2.06%  '
Set-TextValue $ws "E17" '  -1.17%  '
Set-TextValue $ws "D18" '2.610.49'
Set-TextValue $ws "E18" '  +1.03%  '
Set-TextValue $ws "D19" '339.24'
Set-TextValue $ws "E19" '  +1.19%  '
Set-TextValue $ws "D21" '10.29'
Set-TextValue $ws "E21" '  -0.78%  '
Set-TextValue $ws "E22" '  +1.98%  '
Set-TextValue $ws "D23" '0.999'
Set-TextValue $ws "E23" '  +0.01%  '
Set-TextValue $ws "D24" '65.22'
Set-TextValue $ws "E24" '  +1.30%  '
Set-TextValue $ws "D25" '0.167'
Set-TextValue $ws "E25" '  +0.41%  '
Set-TextValue $ws "E26" '  -2.85%  '
Set-TextValue $ws "D27" '2.720.54'
Set-TextValue $ws "E27" '  -0.18%  '
Set-TextValue $ws "E28" '  -0.03%  '
Set-TextValue $ws "D29" '7.02'
Set-TextValue $ws "E29" '  -1.32%  '
Set-TextValue $ws "D30" '0.0₃0749'
Set-TextValue $ws "E30" '  -5.29%  '
Set-TextValue $ws "E31" '  -0.07%  '
Set-TextValue $ws "D32" '6.23'
Set-TextValue $ws "E32" '  -6.32%  '
Set-TextValue $ws "E33" '  +0.32%  '
Set-TextValue $ws "D34" '18.81'
Set-TextValue $ws "E34" '  +0.25%  '
Set-TextValue $ws "D35" '149.85'
Set-TextValue $ws "E35" '  -0.49%  '
Set-TextValue $ws "E36" '  -1.71%  '
Set-TextValue $ws "E37" '  -3.57%  '
Set-TextValue $ws "D38" '0.867'
Set-TextValue $ws "E38" '  -2.57%  '
Set-TextValue $ws "D39" '0.859'
Set-TextValue $ws "E39" '  +1.19%  '
Set-TextValue $ws "B40" 'OKB'
Set-TextValue $ws "C40" 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws "D40" '36.02'
Set-TextValue $ws "E40" '  -0.66%  '
Set-TextValue $ws "B41" 'Stacks'
Set-TextValue $ws "C41" 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws "D41" '1.46'
Set-TextValue $ws "E41" '  +1.63%  '
Set-TextValue $ws "E42" '  -1.94%  '
Set-TextValue $ws "E43" '  -0.37%  '
Set-TextValue $ws "D44" '272.99'
Set-TextValue $ws "E44" '  +1.50%  '
Set-TextValue $ws "E45" '  +0.12%  '
Set-TextValue $ws "D46" '0.0959'
Set-TextValue $ws "E46" '  -0.66%  '
Set-TextValue $ws "E47" '  +0.56%  '
Set-TextValue $ws "D48" '18.80'
Set-TextValue $ws "E48" '  -1.78%  '
Set-TextValue $ws "E49" '  -1.72%  '
Set-TextValue $ws "D50" '18.95'
Set-TextValue $ws "E50" '  +4.08%  '
Set-TextValue $ws "B51" 'Maker'
Set-TextValue $ws "C51" 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws "D51" '1.974.01'
Set-TextValue $ws "E51" '  -3.16%  '
